$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Assesment NO" (column B) result for this rubric changes from 1 to 3.
# Force text storage (so it lands as a shared string, not a number) and
# then drop the temporary number-format so no stray style sticks around.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "3"
$ws.Range("B2").ClearFormats()

# "Component Marks" (column E) result for this rubric changes from 3 to 10
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "10"
$ws.Range("E2").ClearFormats()
